$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Nilai UTS" column (C) is a text column being removed entirely;
# deleting it shifts Praktikum/Tugas (D/E) left into C/D, carrying their
# widths and header text along automatically.
$ws.Range("C1").EntireColumn.Delete()

# Replace the remaining numeric scores with their new values (text -> num).
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 3

$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = 3

$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 3

$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 4
